{"js": "// Replace the two Word field-code paragraphs ({ m:userdoc 'zone1' } and\n// { m:enduserdoc }) with paragraphs containing the same literal text as\n// plain runs (w:t) instead of field characters (w:fldChar/w:instrText).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Grab the Fields collection of every paragraph (kept in a plain array,\n// not attached to the paragraph proxy objects).\nconst fieldCollections = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const fields = paragraphs.items[i].fields;\n  fields.load(\"items\");\n  fieldCollections.push(fields);\n}\nawait context.sync();\n\n// Load the field code for every field found so we can tell the two\n// paragraphs apart.\nfor (let i = 0; i < fieldCollections.length; i++) {\n  const fields = fieldCollections[i];\n  for (let j = 0; j < fields.items.length; j++) {\n    fields.items[j].load(\"code\");\n  }\n}\nawait context.sync();\n\nlet zoneParagraph = null;\nlet endParagraph = null;\n\nfor (let i = 0; i < fieldCollections.length; i++) {\n  const fields = fieldCollections[i];\n  if (fields.items.length === 0) {\n    continue;\n  }\n  const code = fields.items[0].code || \"\";\n  if (code.indexOf(\"enduserdoc\") !== -1) {\n    endParagraph = paragraphs.items[i];\n  } else if (code.indexOf(\"userdoc\") !== -1) {\n    zoneParagraph = paragraphs.items[i];\n  }\n}\n\nif (!zoneParagraph || !endParagraph) {\n  throw new Error(\"Could not locate the expected field paragraphs.\");\n}\n\nconst pkgNs = 'xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"';\nconst wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction wrapOoxml(paragraphXml) {\n  return '<pkg:package ' + pkgNs + '>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document ' + wNs + '><w:body>' + paragraphXml + '<w:sectPr/></w:body></w:document></pkg:xmlData>' +\n    '</pkg:part></pkg:package>';\n}\n\n// { m:userdoc 'zone1' } -> literal text runs: \"{\", \"m\", \":userdoc 'zone1'\", \"}\"\nconst zoneXml = '<w:p>' +\n  '<w:r><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t>:userdoc \\'zone1\\'</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>';\n\n// { m: enduserdoc } -> literal text runs around the preserved bookmark:\n// \"{m:\", <bookmark _GoBack>, \"enduserdoc}\"\nconst endXml = '<w:p>' +\n  '<w:r><w:t>{m:</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  '<w:r><w:t xml:space=\"preserve\">enduserdoc}</w:t></w:r>' +\n  '</w:p>';\n\nzoneParagraph.insertOoxml(wrapOoxml(zoneXml), Word.InsertLocation.replace);\nendParagraph.insertOoxml(wrapOoxml(endXml), Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the two Word field-code paragraphs ({ m:userdoc 'zone1' } and\n# { m:enduserdoc }) with paragraphs containing the same literal text as\n# plain runs (w:t) instead of field characters (w:fldChar/w:instrText).\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndexForPos($doc, $pos) {\n    $idx = 0\n    $foundIdx = -1\n    foreach ($p in $doc.Paragraphs) {\n        $idx = $idx + 1\n        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {\n            $foundIdx = $idx\n        }\n    }\n    return $foundIdx\n}\n\n# Locate the paragraph index (1-based) holding each field, based on the\n# text of its field code (order-independent of how fields are enumerated).\n$zoneParaIndex = -1\n$endParaIndex = -1\n\nforeach ($f in $d.Fields) {\n    $code = $f.Code.Text\n    $pIdx = Find-ParagraphIndexForPos $d $f.Code.Start\n    if ($code -like \"*enduserdoc*\") {\n        $endParaIndex = $pIdx\n    } elseif ($code -like \"*userdoc*\") {\n        $zoneParaIndex = $pIdx\n    }\n}\n\nif ($zoneParaIndex -eq -1 -or $endParaIndex -eq -1) {\n    throw \"Could not locate the expected field paragraphs.\"\n}\n\n$pkgNs = 'xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"'\n$wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\nfunction Wrap-Ooxml($paragraphXml) {\n    return '<pkg:package ' + $pkgNs + '>' + `\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n        '<pkg:xmlData><w:document ' + $wNs + '><w:body>' + $paragraphXml + '<w:sectPr/></w:body></w:document></pkg:xmlData>' + `\n        '</pkg:part></pkg:package>'\n}\n\n# { m:userdoc 'zone1' } -> literal text runs: \"{\", \"m\", \":userdoc 'zone1'\", \"}\"\n$zoneXml = '<w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:userdoc ''zone1''</w:t></w:r><w:r><w:t xml:space=\"preserve\">}</w:t></w:r></w:p>'\n\n# { m: enduserdoc } -> literal text runs around the preserved bookmark:\n# \"{m:\", <bookmark _GoBack>, \"enduserdoc}\"\n$endXml = '<w:p><w:r><w:t>{m:</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t xml:space=\"preserve\">enduserdoc}</w:t></w:r></w:p>'\n\n$zonePara = $d.Paragraphs.Item($zoneParaIndex)\n$zonePara.Range.InsertXML((Wrap-Ooxml $zoneXml))\n\n# Re-fetch the end paragraph by index after the first edit (paragraph\n# count/text didn't change, so the index stays valid).\n$endPara = $d.Paragraphs.Item($endParaIndex)\n$endPara.Range.InsertXML((Wrap-Ooxml $endXml))\n"}
